$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('Z1').Value = 2.025
$ws.Range('Z2').Value = 2.02500000000000001
$ws.Range('Z3').Value = 2.0250000001
